$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.409.33'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.564.85'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.787.68'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '1.567.65'
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.49'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '27.403.29'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '212.95'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.26'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.53'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.11'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.98'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").Value = '1.369.76'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.94'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("D47").Value = '1.700.27'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.53'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("D49").Value = '0.0₇0992'
$ws.Range("E49").Value = '  -2.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0955'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0495'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.82%  '
